$d = $word.ActiveDocument

$rng0 = $d.Content
$rng0.Find.Execute('Zaměstnanec: Jméno a příjmení: [[PERSON_2]] Datum narození: 15.3.1992 Místo narození: Brno Rodné číslo: [[BIRTH_ID_1]] Číslo OP[[PHONE_1]] [[ADDRESS_2]][[PHONE_2]], [[EMAIL_1]] Číslo účtu: [[BANK_1]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng0.Text = 'Zaměstnanec: Jméno a příjmení: [[PERSON_2]] Datum narození: 15.3.1992 Místo narození: Brno Rodné číslo: [[BIRTH_ID_1]] Číslo OP[[PHONE_1]] [[ADDRESS_2]]: +420 [[AMOUNT_1]], [[EMAIL_1]] Číslo účtu: [[BANK_1]]'

$rng1 = $d.Content
$rng1.Find.Execute('Objednatel: [[PERSON_3]] [[ADDRESS_3]]: [[BIRTH_ID_2]] [[PHONE_3]] E-mail: [[EMAIL_2]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng1.Text = 'Objednatel: [[PERSON_3]] [[ADDRESS_3]]: [[BIRTH_ID_2]] Tel.: [[AMOUNT_2]] E-mail: [[EMAIL_2]]'

$rng2 = $d.Content
$rng2.Find.Execute('Zhotovitel: [[PERSON_4]] [[ADDRESS_4]] [[ICO_2]] DIČ: [[DIC_2]] Datum narození: 23.09.1985 Bankovní spojení: [[BANK_2]] [[PHONE_4]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng2.Text = 'Zhotovitel: [[PERSON_4]] [[ADDRESS_4]] [[ICO_2]] DIČ: [[DIC_2]] Datum narození: 23.09.1985 Bankovní spojení: [[BANK_2]] Telefon: +420 [[AMOUNT_3]]'

$rng3 = $d.Content
$rng3.Find.Execute('[[PERSON_5]] se zavazuje pro pana doktora Krajíčka provést kompletní renovaci ordinace. MUDr. Krajíček uhradí [[PERSON_4]] částku 285 000 Kč.', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng3.Text = '[[PERSON_5]] se zavazuje pro pana doktora Krajíčka provést kompletní renovaci ordinace. MUDr. Krajíček uhradí [[PERSON_4]] částku [[AMOUNT_4]].'

$rng4 = $d.Content
$rng4.Find.Execute('Pronajímatel: [[PERSON_6]], nar. 5.8.1968 (r.č. [[BIRTH_ID_3]]) Trvale bytem: [[ADDRESS_5]]: [[ID_CARD_1]], vydán 12.2.2020 Kontakt: [[EMAIL_3]], [[PHONE_5]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng4.Text = 'Pronajímatel: [[PERSON_6]], nar. 5.8.1968 (r.č. [[BIRTH_ID_3]]) Trvale bytem: [[ADDRESS_5]]: [[ID_CARD_1]], vydán 12.2.2020 Kontakt: [[EMAIL_3]], tel. [[AMOUNT_5]]'

$rng5 = $d.Content
$rng5.Find.Execute('Nájemce: [[PERSON_7]] (r.č. [[BIRTH_ID_4]]) Bytem: [[ADDRESS_6]]: [[BANK_3]] [[PHONE_6]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng5.Text = 'Nájemce: [[PERSON_7]] (r.č. [[BIRTH_ID_4]]) Bytem: [[ADDRESS_6]]: [[BANK_3]] Mobil: +420 [[AMOUNT_6]]'

$rng6 = $d.Content
$rng6.Find.Execute('[[PERSON_8]] pronajímá panu Procházkovi byt 2+kk v ulici Mánesova 87, Brno. Procházka se zavazuje platit Horváthové měsíční nájemné 15 000 Kč.', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng6.Text = '[[PERSON_8]] pronajímá panu Procházkovi byt 2+kk v ulici Mánesova 87, Brno. Procházka se zavazuje platit Horváthové měsíční nájemné [[AMOUNT_7]].'

$rng7 = $d.Content
$rng7.Find.Execute('Zprostředkovatel: RealEstate Pro s.r.o. [[ADDRESS_7]] [[ICO_3]] Zastoupená: [[PERSON_9]], Ph.D., jednatelem E-mail: [[EMAIL_4]] [[PHONE_7]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng7.Text = 'Zprostředkovatel: RealEstate Pro s.r.o. [[ADDRESS_7]] [[ICO_3]] Zastoupená: [[PERSON_9]], Ph.D., jednatelem E-mail: [[EMAIL_4]] Tel.: +420 [[AMOUNT_8]]'

$rng8 = $d.Content
$rng8.Find.Execute('Klient: Mgr. Kateřina Malá, Ph.D. Narozena: 12.11.1985 (RČ: [[BIRTH_ID_5]]) Bytem: [[ADDRESS_8]]: [[ID_CARD_2]] E-mail: [[EMAIL_5]] [[PHONE_8]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng8.Text = 'Klient: Mgr. Kateřina Malá, Ph.D. Narozena: 12.11.1985 (RČ: [[BIRTH_ID_5]]) Bytem: [[ADDRESS_8]]: [[ID_CARD_2]] E-mail: [[EMAIL_5]] Telefon: [[AMOUNT_9]]'

$rng9 = $d.Content
$rng9.Find.Execute('Prodávající: Martin "Marty" Král Bytem: [[ADDRESS_9]] Nar.: 4.6.1979, RČ: [[BIRTH_ID_6]] OP: [[ID_CARD_3]] [[PHONE_9]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng9.Text = 'Prodávající: Martin "Marty" Král Bytem: [[ADDRESS_9]] Nar.: 4.6.1979, RČ: [[BIRTH_ID_6]] OP: [[ID_CARD_3]] Tel.: [[AMOUNT_10]]'

$rng10 = $d.Content
$rng10.Find.Execute('Kupující: [[PERSON_12]] (rozená Nová) [[ADDRESS_10]]: [[BIRTH_ID_7]] Datum narození: 23.4.1988 Kontakt: [[EMAIL_6]],[[PHONE_10]] Bankovní účet: [[BANK_4]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng10.Text = 'Kupující: [[PERSON_12]] (rozená Nová) [[ADDRESS_10]]: [[BIRTH_ID_7]] Datum narození: 23.4.1988 Kontakt: [[EMAIL_6]], [[AMOUNT_11]] Bankovní účet: [[BANK_4]]'

$rng11 = $d.Content
$rng11.Find.Execute('[[PERSON_13]] (dále jen "Marty") prodává paní Beránková (dříve [[PERSON_14]]) osobní automobil Škoda Octavia. Petra uhradí Martinovi částku 250 000 Kč na účet [[BANK_5]].', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng11.Text = '[[PERSON_13]] (dále jen "Marty") prodává paní Beránková (dříve [[PERSON_14]]) osobní automobil Škoda Octavia. Petra uhradí Martinovi částku [[AMOUNT_12]] na účet [[BANK_5]].'

$rng12 = $d.Content
$rng12.Find.Execute('Strana A: [[PERSON_15]] (samostatný konzultant) [[ICO_4]] [[ADDRESS_11]] Rodné číslo: [[BIRTH_ID_8]] E-mail: [[EMAIL_7]] [[PHONE_11]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng12.Text = 'Strana A: [[PERSON_15]] (samostatný konzultant) [[ICO_4]] [[ADDRESS_11]] Rodné číslo: [[BIRTH_ID_8]] E-mail: [[EMAIL_7]] Tel.: +420 [[AMOUNT_13]]'

$rng13 = $d.Content
$rng13.Find.Execute('Strana B: Innovation Labs a.s. [[ADDRESS_12]] [[ICO_5]] Zastoupena: [[PERSON_16]], MBA (osobní kontakt: [[EMAIL_8]],[[PHONE_12]])', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng13.Text = 'Strana B: Innovation Labs a.s. [[ADDRESS_12]] [[ICO_5]] Zastoupena: [[PERSON_16]], MBA (osobní kontakt: [[EMAIL_8]], [[AMOUNT_14]])'

$rng14 = $d.Content
$rng14.Find.Execute('Případné spory budou řešit mediátorkou [[PERSON_19]] (nar. 15.5.1975, kontakt[[PHONE_13]]).', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng14.Text = 'Případné spory budou řešit mediátorkou [[PERSON_19]] (nar. 15.5.1975, kontakt: [[AMOUNT_15]]).'

$rng15 = $d.Content
$rng15.Find.Execute('Zaměstnanec: Nguyễn Thị Lan Bytem: [[ADDRESS_14]] Rodné číslo: [[BIRTH_ID_9]] Číslo OP: [[ID_CARD_4]] Narozena: 12.3.1996 Kontakt: [[EMAIL_10]] [[PHONE_14]] Číslo účtu: [[BANK_6]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng15.Text = 'Zaměstnanec: Nguyễn Thị Lan Bytem: [[ADDRESS_14]] Rodné číslo: [[BIRTH_ID_9]] Číslo OP: [[ID_CARD_4]] Narozena: 12.3.1996 Kontakt: [[EMAIL_10]] Telefon: +420 [[AMOUNT_16]] Číslo účtu: [[BANK_6]]'

$rng16 = $d.Content
$rng16.Find.Execute('Věřitel: [[PERSON_21]] Bytem: [[ADDRESS_15]]: [[BIRTH_ID_10]] OP: [[ID_CARD_5]] [[PHONE_15]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng16.Text = 'Věřitel: [[PERSON_21]] Bytem: [[ADDRESS_15]]: [[BIRTH_ID_10]] OP: [[ID_CARD_5]] Tel.: [[AMOUNT_17]]'

$rng17 = $d.Content
$rng17.Find.Execute('[[PERSON_22]] [[ADDRESS_16]]: [[BIRTH_ID_11]] Občanský průkaz: [[ID_CARD_6]] [[PHONE_16]] E-mail: [[EMAIL_11]]', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng17.Text = '[[PERSON_22]] [[ADDRESS_16]]: [[BIRTH_ID_11]] Občanský průkaz: [[ID_CARD_6]] Telefon: [[AMOUNT_18]] E-mail: [[EMAIL_11]]'

$rng18 = $d.Content
$rng18.Find.Execute('[[PERSON_23]] poskytuje panu Novému úvěr ve výši 150 000 Kč. Nový splácí Novotnému v měsíčních splátkách. V případě prodlení zaplatí [[PERSON_22]] [[PERSON_21]] penále.', $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng18.Text = '[[PERSON_23]] poskytuje panu Novému úvěr ve výši [[AMOUNT_19]]. Nový splácí Novotnému v měsíčních splátkách. V případě prodlení zaplatí [[PERSON_22]] [[PERSON_21]] penále.'

